# Applies the numeric updates from the diff (Case_3_121, 380 kV case) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 10.39602535399267
$ws.Range("C2").Value = 4.135582061848216
$ws.Range("D2").Value = 10.05864041608433
$ws.Range("F2").Value = 34.79411174578443
$ws.Range("G2").Value = 3.665364992533431
$ws.Range("J2").Value = 11.33038245765799
$ws.Range("K2").Value = 9.7921287399933
$ws.Range("M2").Value = 16.32970447612944
$ws.Range("N2").Value = 20.07931252617129
$ws.Range("O2").Value = 25.86532493136859
# Row 3
$ws.Range("B3").Value = 10.16711443461404
$ws.Range("C3").Value = 3.97644282478712
$ws.Range("D3").Value = 10.01775615280918
$ws.Range("F3").Value = 34.84360220454653
$ws.Range("G3").Value = 3.667143155626547
$ws.Range("J3").Value = 11.34755289905068
$ws.Range("K3").Value = 9.643952018225587
$ws.Range("M3").Value = 16.26374454266087
$ws.Range("N3").Value = 20.14002625001507
$ws.Range("O3").Value = 25.93428203033342
# Row 4
$ws.Range("B4").Value = 10.02598674484173
$ws.Range("C4").Value = 3.874652717405173
$ws.Range("D4").Value = 9.994364236923495
$ws.Range("F4").Value = 34.88138645691993
$ws.Range("G4").Value = 3.668293479806635
$ws.Range("J4").Value = 11.35968546778925
$ws.Range("K4").Value = 9.553513760216759
$ws.Range("M4").Value = 16.22591919658767
$ws.Range("N4").Value = 20.17904413660762
$ws.Range("O4").Value = 25.98143678338791
# Row 5
$ws.Range("B5").Value = 9.968417188751427
$ws.Range("C5").Value = 3.832184830550421
$ws.Range("D5").Value = 9.985269387815176
$ws.Range("F5").Value = 34.8986418649825
$ws.Range("G5").Value = 3.668777008248083
$ws.Range("J5").Value = 11.3650295072906
$ws.Range("K5").Value = 9.516841130804
$ws.Range("M5").Value = 16.21118917799142
$ws.Range("N5").Value = 20.19538291572005
$ws.Range("O5").Value = 26.00186190490783
# Row 6
$ws.Range("B6").Value = 9.958856809948324
$ws.Range("C6").Value = 3.825074624132272
$ws.Range("D6").Value = 9.983785823834895
$ws.Range("F6").Value = 34.90161926303517
$ws.Range("G6").Value = 3.66885819070324
$ws.Range("J6").Value = 11.36594103829449
$ws.Range("K6").Value = 9.51076396408231
$ws.Range("M6").Value = 16.20878491347097
$ws.Range("N6").Value = 20.1981224901104
$ws.Range("O6").Value = 26.00532646488492
# Row 7
$ws.Range("B7").Value = 10.02521046477322
$ws.Range("C7").Value = 3.87408392577054
$ws.Range("D7").Value = 9.994239799786277
$ws.Range("F7").Value = 34.8816116494013
$ws.Range("G7").Value = 3.668299941008928
$ws.Range("J7").Value = 11.35975591996589
$ws.Range("K7").Value = 9.553018384895212
$ws.Range("M7").Value = 16.22571775764457
$ws.Range("N7").Value = 20.17926270922635
$ws.Range("O7").Value = 25.98170734994255
# Row 8
$ws.Range("B8").Value = 10.31726695592237
$ws.Range("C8").Value = 4.081578463330859
$ws.Range("D8").Value = 10.04419286702543
$ws.Range("F8").Value = 34.80963978441845
$ws.Range("G8").Value = 3.665965983734516
$ws.Range("J8").Value = 11.33597294840288
$ws.Range("K8").Value = 9.740951809745644
$ws.Range("M8").Value = 16.30641364586157
$ws.Range("N8").Value = 20.09988637902383
$ws.Range("O8").Value = 25.88810112286009
# Row 9
$ws.Range("B9").Value = 10.8818010481764
$ws.Range("C9").Value = 4.454691219132172
$ws.Range("D9").Value = 10.15538753308379
$ws.Range("F9").Value = 34.72726346859058
$ws.Range("G9").Value = 3.661851395283432
$ws.Range("J9").Value = 11.30194227954211
$ws.Range("K9").Value = 10.11182333382776
$ws.Range("M9").Value = 16.48533302064598
$ws.Range("N9").Value = 19.95797147189642
$ws.Range("O9").Value = 25.74280919488118
# Row 10
$ws.Range("B10").Value = 11.28697372016145
$ws.Range("C10").Value = 4.706609648543788
$ws.Range("D10").Value = 10.24466535458345
$ws.Range("F10").Value = 34.70262792124444
$ws.Range("G10").Value = 3.659107331724448
$ws.Range("J10").Value = 11.28461491602058
$ws.Range("K10").Value = 10.38312575776694
$ws.Range("M10").Value = 16.62861397317426
$ws.Range("N10").Value = 19.86199897744248
$ws.Range("O10").Value = 25.65947929521219
# Row 11
$ws.Range("B11").Value = 11.46832378864206
$ws.Range("C11").Value = 4.816125807520686
$ws.Range("D11").Value = 10.28681421994195
$ws.Range("F11").Value = 34.6992143821803
$ws.Range("G11").Value = 3.657918934975835
$ws.Range("J11").Value = 11.27839549108628
$ws.Range("K11").Value = 10.50576083118627
$ws.Range("M11").Value = 16.69618603384941
$ws.Range("N11").Value = 19.82012106566284
$ws.Range("O11").Value = 25.62667110521156
# Row 12
$ws.Range("B12").Value = 11.53650189730244
$ws.Range("C12").Value = 4.856847798498057
$ws.Range("D12").Value = 10.30298567148032
$ws.Range("F12").Value = 34.69904135237303
$ws.Range("G12").Value = 3.657477485692902
$ws.Range("J12").Value = 11.27627909273758
$ws.Range("K12").Value = 10.55204454432226
$ws.Range("M12").Value = 16.72210190702038
$ws.Range("N12").Value = 19.80451767950189
$ws.Range("O12").Value = 25.61498184686521
# Row 13
$ws.Range("B13").Value = 11.52184159265686
$ws.Range("C13").Value = 4.84811118417208
$ws.Range("D13").Value = 10.29949366105405
$ws.Range("F13").Value = 34.69902884689103
$ws.Range("G13").Value = 3.657572179164179
$ws.Range("J13").Value = 11.27672428425563
$ws.Range("K13").Value = 10.54208410447951
$ws.Range("M13").Value = 16.71650614567154
$ws.Range("N13").Value = 19.80786683078974
$ws.Range("O13").Value = 25.61746665253674
# Row 14
$ws.Range("B14").Value = 11.47394314518577
$ws.Range("C14").Value = 4.819491145464451
$ws.Range("D14").Value = 10.28814049446075
$ws.Range("F14").Value = 34.69917771688472
$ws.Range("G14").Value = 3.657882445143219
$ws.Range("J14").Value = 11.27821659080953
$ws.Range("K14").Value = 10.50957197817766
$ws.Range("M14").Value = 16.69831167379665
$ws.Range("N14").Value = 19.81883226430754
$ws.Range("O14").Value = 25.62569469587718
# Row 15
$ws.Range("B15").Value = 11.44453750781591
$ws.Range("C15").Value = 4.801862449629326
$ws.Range("D15").Value = 10.28121346666583
$ws.Range("F15").Value = 34.6994146663917
$ws.Range("G15").Value = 3.658073606954351
$ws.Range("J15").Value = 11.27916175404393
$ws.Range("K15").Value = 10.48963585641559
$ws.Range("M15").Value = 16.687209231719
$ws.Range("N15").Value = 19.82558206640476
$ws.Range("O15").Value = 25.63083029744224
# Row 16
$ws.Range("B16").Value = 11.27505640775234
$ws.Range("C16").Value = 4.699348584998692
$ws.Range("D16").Value = 10.24194086347084
$ws.Range("F16").Value = 34.70300774958741
$ws.Range("G16").Value = 3.659186197857772
$ws.Range("J16").Value = 11.28505480225415
$ws.Range("K16").Value = 10.37509181390562
$ws.Range("M16").Value = 16.62424474958678
$ws.Range("N16").Value = 19.86477152563385
$ws.Range("O16").Value = 25.66172610353114
# Row 17
$ws.Range("B17").Value = 11.17027532637172
$ws.Range("C17").Value = 4.635143750621423
$ws.Range("D17").Value = 10.21823458190106
$ws.Range("F17").Value = 34.70720749680923
$ws.Range("G17").Value = 3.659884046318761
$ws.Range("J17").Value = 11.28909566413928
$ws.Range("K17").Value = 10.30459134323436
$ws.Range("M17").Value = 16.58621951964303
$ws.Range("N17").Value = 19.88926818045543
$ws.Range("O17").Value = 25.68198659795639
# Row 18
$ws.Range("B18").Value = 11.10973275700681
$ws.Range("C18").Value = 4.5977373252529
$ws.Range("D18").Value = 10.20474456301401
$ws.Range("F18").Value = 34.71035672813804
$ws.Range("G18").Value = 3.660291070031744
$ws.Range("J18").Value = 11.29157640414303
$ws.Range("K18").Value = 10.26397023612944
$ws.Range("M18").Value = 16.56457464748382
$ws.Range("N18").Value = 19.90352566195465
$ws.Range("O18").Value = 25.69411977371827
# Row 19
$ws.Range("B19").Value = 11.08918898290078
$ws.Range("C19").Value = 4.584990751121203
$ws.Range("D19").Value = 10.20020232138622
$ws.Range("F19").Value = 34.71154902986176
$ws.Range("G19").Value = 3.660429851169302
$ws.Range("J19").Value = 11.29244323644765
$ws.Range("K19").Value = 10.25020575561714
$ws.Range("M19").Value = 16.55728540942096
$ws.Range("N19").Value = 19.90838183032692
$ws.Range("O19").Value = 25.69831023687434
# Row 20
$ws.Range("B20").Value = 11.18145845397043
$ws.Range("C20").Value = 4.64202802751934
$ws.Range("D20").Value = 10.22074320100233
$ws.Range("F20").Value = 34.70668450229847
$ws.Range("G20").Value = 3.659809175769005
$ws.Range("J20").Value = 11.28864930827036
$ws.Range("K20").Value = 10.31210391484706
$ws.Range("M20").Value = 16.5902440722284
$ws.Range("N20").Value = 19.88664312662453
$ws.Range("O20").Value = 25.67978015885109
# Row 21
$ws.Range("B21").Value = 11.48802604567377
$ws.Range("C21").Value = 4.827918016827801
$ws.Range("D21").Value = 10.29146956283606
$ws.Range("F21").Value = 34.69910361626456
$ws.Range("G21").Value = 3.657791080225322
$ws.Range("J21").Value = 11.27777178738071
$ws.Range("K21").Value = 10.51912613411864
$ws.Range("M21").Value = 16.70364706710393
$ws.Range("N21").Value = 19.81560454274564
$ws.Range("O21").Value = 25.62325797687098
# Row 22
$ws.Range("B22").Value = 11.68547056926913
$ws.Range("C22").Value = 4.945033095779664
$ws.Range("D22").Value = 10.33891582411508
$ws.Range("F22").Value = 34.70067400252727
$ws.Range("G22").Value = 3.656522076120742
$ws.Range("J22").Value = 11.27205418437173
$ws.Range("K22").Value = 10.65350179172834
$ws.Range("M22").Value = 16.77966525559049
$ws.Range("N22").Value = 19.77066177035132
$ws.Range("O22").Value = 25.59059912371918
# Row 23
$ws.Range("B23").Value = 11.58037895748449
$ws.Range("C23").Value = 4.882932275575184
$ws.Range("D23").Value = 10.31348448081055
$ws.Range("F23").Value = 34.69923935205707
$ws.Range("G23").Value = 3.657194811803758
$ws.Range("J23").Value = 11.27497858881323
$ws.Range("K23").Value = 10.58188144956507
$ws.Range("M23").Value = 16.73892432073007
$ws.Range("N23").Value = 19.79451307142094
$ws.Range("O23").Value = 25.6076376164378
# Row 24
$ws.Range("B24").Value = 11.17640349900157
$ws.Range("C24").Value = 4.638917183881428
$ws.Range("D24").Value = 10.21960862044923
$ws.Range("F24").Value = 34.70691865936683
$ws.Range("G24").Value = 3.659843006608078
$ws.Range("J24").Value = 11.28885061484023
$ws.Range("K24").Value = 10.30870775652175
$ws.Range("M24").Value = 16.58842389689858
$ws.Range("N24").Value = 19.88782937127143
$ws.Range("O24").Value = 25.68077617884781
# Row 25
$ws.Range("B25").Value = 10.73044655665825
$ws.Range("C25").Value = 4.357552632141795
$ws.Range("D25").Value = 10.12393777314282
$ws.Range("F25").Value = 34.74324865467417
$ws.Range("G25").Value = 3.662915309994534
$ws.Range("J25").Value = 11.3097996683236
$ws.Range("K25").Value = 10.0115146820389
$ws.Range("M25").Value = 16.43479449563146
$ws.Range("N25").Value = 19.9949008822114
$ws.Range("O25").Value = 25.7780081464716

Write-Output "Applied 240 cell updates (B/C/D/F/G/J/K/M/N/O across rows 2-25)."
